$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of rows 2 and 3 for the columns that differ between
# the two records (A, B, E, F, G, H, P, Q, R, S, AC). The remaining
# columns already hold identical values in both rows, so only these
# need to be exchanged.
$cols = @("A","B","E","F","G","H","P","Q","R","S","AC")

foreach ($col in $cols) {
    $cell2 = $ws.Range("$col`2")
    $cell3 = $ws.Range("$col`3")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $tmp
}
